# Apply updated "F" column (view/heat count) values across the workbook,
# matching the regenerated data snapshot described in the commit
# "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 285
$ws1.Range("F3").Value = 971
$ws1.Range("F6").Value = 3198
$ws1.Range("F16").Value = 1335
$ws1.Range("F17").Value = 1335
$ws1.Range("F23").Value = 444
$ws1.Range("F26").Value = 334
$ws1.Range("F29").Value = 169
$ws1.Range("F31").Value = 396
$ws1.Range("F33").Value = 5204

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F17").Value = 54
$ws2.Range("F34").Value = 35
$ws2.Range("F35").Value = 35

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 285
$ws4.Range("F6").Value = 971
$ws4.Range("F24").Value = 1335
$ws4.Range("F25").Value = 1335
$ws4.Range("F31").Value = 444
$ws4.Range("F32").Value = 54
$ws4.Range("F35").Value = 334
$ws4.Range("F39").Value = 169
$ws4.Range("F42").Value = 396
$ws4.Range("F43").Value = 5204
$ws4.Range("F47").Value = 35
